$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Date" column (BF) contains an incorrectly formatted/off-by-one-day
# value "6-9-2011-12" for every data row (rows 2-31). Replace it with the
# corrected date string "2012-06-09" to match the rest of the training data.
#
# Note: assigning the literal text "2012-06-09" straight into .Value /
# .Value2 causes Excel to auto-recognize it as a date and silently store
# it as a date serial number instead of the plain text string. To keep it
# as literal text (matching the original inline-string cell type) we first
# write it as a text-literal formula and then convert the cell to a static
# value via copy / paste-special values, which avoids the date-recognition
# that a direct value assignment would trigger.
$rng = $ws.Range("BF2:BF31")
$rng.Formula = '="2012-06-09"'
$rng.Copy()
$rng.PasteSpecial(-4163) # xlPasteValues
$excel.CutCopyMode = 0
